{"js": "// Insert a new bulleted list item \"Check for valid moves\" right after the\n// existing \"Be able to update/remove where the player has cells.\" item,\n// matching the author's commit (end of 2nd iteration step list).\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst anchorText = \"Be able to update/remove where the player has cells.\";\nlet anchor = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const t = paragraphs.items[i].text.trim();\n  if (t === anchorText) {\n    anchor = paragraphs.items[i];\n    break;\n  }\n}\n\nif (!anchor) {\n  throw new Error(\"Could not find anchor paragraph: \" + anchorText);\n}\n\n// insertParagraph inherits the anchor paragraph's formatting (style +\n// numbering), so the new paragraph keeps the ListParagraph style / bullet\n// list (numId 1, ilvl 0) just like the diff shows.\nconst newPara = anchor.insertParagraph(\"Check for valid moves\", Word.InsertLocation.after);\n\nawait context.sync();\n", "ps1": "# Insert a new bulleted list item \"Check for valid moves\" right after the\n# existing \"Be able to update/remove where the player has cells.\" item\n# (end of the 2nd-iteration \"Steps involved\" list), matching the author's\n# commit (\"Finished 2nd iteration...\").\n\n$d = $word.ActiveDocument\n\n$anchorText = \"Be able to update/remove where the player has cells.\"\n\n$anchorIndex = -1\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($t -eq $anchorText) {\n        $anchorIndex = $i\n        break\n    }\n}\n\nif ($anchorIndex -eq -1) {\n    throw \"Could not find anchor paragraph: $anchorText\"\n}\n\n$anchorPara = $d.Paragraphs.Item($anchorIndex)\n\n# InsertParagraphAfter inherits the source paragraph's formatting (style +\n# list numbering), so the new paragraph keeps the ListParagraph style /\n# bullet list (numId 1, ilvl 0), matching the diff. It returns nothing (like\n# real Word VBA), so grab the freshly created paragraph by index afterwards.\n$anchorPara.Range.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item($anchorIndex + 1)\n$newPara.Range.Text = \"Check for valid moves\"\n"}
